# Updates 合肥-漫展信息.xlsx to add the "合肥·微光mini动漫派对（免费展）" event
# (2024-08-24) into the "展览" and "全部类型" sheets, refresh a few
# "想去人数" (interest count) figures, and renumber the trailing rows.

function Set-TextValue {
    # Force a literal-text write so date-looking strings (e.g. "2024-08-24")
    # are not auto-coerced into Excel date serials: build the value via a
    # formula returning a string literal, then flatten it to a static value
    # with Paste-Special-Values (xlPasteValues = -4163). This also avoids
    # creating any new cell style (no NumberFormat / quote-prefix touched).
    param($ws, $addr, [string]$val)
    $escaped = $val -replace '"', '""'
    $cell = $ws.Range($addr)
    $cell.Formula = "=""$escaped"""
    $cell.Copy()
    $cell.PasteSpecial(-4163)
}

function Insert-EventRow {
    param($ws, [int]$rowNum, [int]$seq)

    # Push existing row $rowNum (and below) down by one.
    $ws.Rows.Item($rowNum).Insert()

    # Clone the formatting of the row above (same look as every other data
    # row: bold/bordered/centered column A, default formatting elsewhere).
    $srcRange = "A" + ($rowNum - 1) + ":I" + ($rowNum - 1)
    $dstRange = "A" + $rowNum + ":I" + $rowNum
    $ws.Range($srcRange).Copy()
    $ws.Range($dstRange).PasteSpecial(-4122)

    $aAddr = "A" + $rowNum
    $bAddr = "B" + $rowNum
    $cAddr = "C" + $rowNum
    $dAddr = "D" + $rowNum
    $eAddr = "E" + $rowNum
    $fAddr = "F" + $rowNum
    $gAddr = "G" + $rowNum
    $hAddr = "H" + $rowNum
    $iAddr = "I" + $rowNum

    $ws.Range($aAddr).Value = $seq
    Set-TextValue $ws $bAddr "2024-08-24"
    Set-TextValue $ws $cAddr "合肥·微光mini动漫派对（免费展）"
    Set-TextValue $ws $dAddr "山林路与山水路交叉路口往东北约70米 伟星星悦广场(肥东店)"
    Set-TextValue $ws $eAddr "2024.08.24 13:00-08.25 19:00"
    $ws.Range($fAddr).Value = 2
    $ws.Range($gAddr).Value = 58
    Set-TextValue $ws $hAddr "https://show.bilibili.com/platform/detail.html?id=90625"
    Set-TextValue $ws $iAddr "//i0.hdslb.com/bfs/openplatform/202408/t7kq4X7h1723471019389.jpeg"
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "展览": refresh counts, insert the new event as row 7 (shifting
# the two rows below it down to 8/9), then renumber their A column.
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")

$ws1.Range("F2").Value = 5369
$ws1.Range("F3").Value = 379
$ws1.Range("F6").Value = 809

Insert-EventRow $ws1 7 6

$ws1.Range("A8").Value = 7
$ws1.Range("A9").Value = 8
$ws1.Range("F8").Value = 325

# ---------------------------------------------------------------------
# Sheet "全部类型": same refresh, insert the new event as row 7 (shifting
# the five rows below it down to 8-12), then renumber their A column.
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")

$ws4.Range("F2").Value = 5369
$ws4.Range("F3").Value = 379
$ws4.Range("F6").Value = 809

Insert-EventRow $ws4 7 6

$ws4.Range("A8").Value = 7
$ws4.Range("A9").Value = 8
$ws4.Range("A10").Value = 9
$ws4.Range("A11").Value = 10
$ws4.Range("A12").Value = 11
$ws4.Range("F9").Value = 325

Write-Output "edit complete"
